$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab) from "SCD0265" to "SCD0016"
$ws.Name = "SCD0016"

# 2. Update TC_ID cells (column B, rows 2-4) from "DGS-280" to "SCD0016-039"
#    (the old "DGS-280" shared string is dropped and a new "SCD0016-039"
#    shared string is appended, which also renumbers the other shared
#    strings used by columns D/E - that happens automatically)
$ws.Range("B2").Value = "SCD0016-039"
$ws.Range("B3").Value = "SCD0016-039"
$ws.Range("B4").Value = "SCD0016-039"

# 3. Widen column B to fit the longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.43

# 4. Update the view: scroll so row 3 is at the top and select B6
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("B6").Select() | Out-Null
